# openfisca_france/assets/grilles_fonction_publique/test_grid.xlsx
#
# 1) Fix the misspelling "redaceur chef" -> "redacteur chef" for the
#    "chef" rows (F23:F29 and F37:F46).
# 2) Rows F30:F36 actually describe the "redacteur" grade (not "redacteur
#    chef"), so fix those to plain "redacteur".
# 3) Update the sheet selection left over from the author's last edit:
#    was L30:L46 (active cell L30) -> now A30:M36 (active cell A30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- F23:F29 -> "redacteur chef" (spelling fix) ---------------------------
foreach ($r in 23..29) {
    $ws.Cells.Item($r, 6).Value = "redacteur chef"
}

# --- F30:F36 -> "redacteur" -------------------------------------------------
foreach ($r in 30..36) {
    $ws.Cells.Item($r, 6).Value = "redacteur"
}

# --- F37:F46 -> "redacteur chef" (spelling fix) ----------------------------
foreach ($r in 37..46) {
    $ws.Cells.Item($r, 6).Value = "redacteur chef"
}

# --- Update the sheet's saved selection ------------------------------------
$ws.Select() | Out-Null
$ws.Range("A30:M36").Select() | Out-Null
